$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# --- Update row 43: mark fm38 test as complete (was "in progress") ---
$ws.Range("H43").Value = "complete"
$ws.Range("I43").Value = "complete"

# --- Add new row 57: test case fm53 for step policies ---
# Copy formatting from the row above (row 56) so the new row matches the
# existing table styling (font/border style index 5 on B,C,H,I).
$ws.Range("B56:C56").Copy()
$ws.Range("B57").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H56:I56").Copy()
$ws.Range("H57").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B57").Value = "fm53"
$ws.Range("C57").Value = "JP Flood step policies with extra expense and debris removal "
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 27
$ws.Range("F57").Value = 1
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = "in progress"
$ws.Range("I57").Value = "in progress"

# --- Update the selection to reflect where the user ended up ---
$ws.Activate()
$ws.Range("H56:I56").Select()
